$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.231460213661194
$ws.Range("B1").Value = 2.530009269714355
$ws.Range("C1").Value = 7.638636112213135
$ws.Range("D1").Value = 2.191313505172729
$ws.Range("E1").Value = 1.150243997573853
